$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coin/Price/Volume table refresh (GitHub Actions data pull).
# A leading apostrophe forces Price-column values that look numeric
# (e.g. "1.70") to stay plain text, matching the original inlineStr cells
# instead of being auto-coerced into a Number by Excel.

$ws.Range('D2').Value = '65.334.58'
$ws.Range('E2').Value = '  -0.58%  '

$ws.Range('D3').Value = '3.543.48'
$ws.Range('E3').Value = '  +2.78%  '

$ws.Range('D4').Value = '''0.999'
$ws.Range('E4').Value = '  -0.06%  '

$ws.Range('D5').Value = '''604.05'
$ws.Range('E5').Value = '  +1.75%  '

$ws.Range('D6').Value = '''140.18'
$ws.Range('E6').Value = '  +2.50%  '

$ws.Range('D7').Value = '3.542.68'
$ws.Range('E7').Value = '  +2.79%  '

$ws.Range('E8').Value = '  +0.03%  '

$ws.Range('D9').Value = '''0.492'
$ws.Range('E9').Value = '  -1.73%  '

$ws.Range('E10').Value = '  +2.40%  '

$ws.Range('D11').Value = '''7.01'
$ws.Range('E11').Value = '  -5.09%  '

$ws.Range('D12').Value = '''0.393'
$ws.Range('E12').Value = '  +3.52%  '

$ws.Range('D13').Value = '4.145.26'
$ws.Range('E13').Value = '  +2.86%  '

$ws.Range('E14').Value = '  +3.03%  '

$ws.Range('D15').Value = '''27.31'
$ws.Range('E15').Value = '  +2.36%  '

$ws.Range('D16').Value = '3.545.19'
$ws.Range('E16').Value = '  +1.79%  '

$ws.Range('E17').Value = '  +1.77%  '

$ws.Range('D18').Value = '65.392.23'
$ws.Range('E18').Value = '  -0.38%  '

$ws.Range('D19').Value = '''10.39'
$ws.Range('E19').Value = '  +4.93%  '

$ws.Range('D20').Value = '''5.97'
$ws.Range('E20').Value = '  +1.91%  '

$ws.Range('D21').Value = '''14.36'
$ws.Range('E21').Value = '  +4.52%  '

$ws.Range('D22').Value = '''396.39'
$ws.Range('E22').Value = '  +0.78%  '

$ws.Range('D23').Value = '''0.575'
$ws.Range('E23').Value = '  +4.14%  '

$ws.Range('D24').Value = '3.683.63'
$ws.Range('E24').Value = '  +2.59%  '

$ws.Range('D25').Value = '''74.08'
$ws.Range('E25').Value = '  +0.72%  '

$ws.Range('E26').Value = '  +0.00%  '

$ws.Range('E27').Value = '  +9.96%  '

$ws.Range('D28').Value = '''7.88'
$ws.Range('E28').Value = '  +9.15%  '

$ws.Range('E29').Value = '  +0.17%  '

$ws.Range('E30').Value = '  +2.28%  '

$ws.Range('D31').Value = '''8.34'
$ws.Range('E31').Value = '  +1.28%  '

$ws.Range('D32').Value = '3.553.58'
$ws.Range('E32').Value = '  +2.86%  '

$ws.Range('E33').Value = '  +0.03%  '

$ws.Range('B34').Value = 'EthereumClassic'
$ws.Range('C34').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D34').Value = '''23.85'
$ws.Range('E34').Value = '  +3.43%  '

$ws.Range('B35').Value = 'Kaspa'
$ws.Range('C35').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D35').Value = '''0.147'
$ws.Range('E35').Value = '  -0.27%  '

$ws.Range('E36').Value = '  +8.82%  '

$ws.Range('D37').Value = '''7.04'
$ws.Range('E37').Value = '  +0.86%  '

$ws.Range('B38').Value = 'ImmutableX'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D38').Value = '''1.57'
$ws.Range('E38').Value = '  +3.42%  '

$ws.Range('B39').Value = 'Monero'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D39').Value = '''169.64'
$ws.Range('E39').Value = '  -1.28%  '

$ws.Range('E40').Value = '  +3.37%  '

$ws.Range('D41').Value = '''0.0818'
$ws.Range('E41').Value = '  +6.07%  '

$ws.Range('E42').Value = '  +0.78%  '

$ws.Range('D43').Value = '''26.56'
$ws.Range('E43').Value = '  +15.46%  '

$ws.Range('D44').Value = '''42.94'
$ws.Range('E44').Value = '  -1.78%  '

$ws.Range('E45').Value = '  -0.07%  '

$ws.Range('D46').Value = '''4.46'
$ws.Range('E46').Value = '  +0.81%  '

$ws.Range('E47').Value = '  +10.11%  '

$ws.Range('D48').Value = '''1.70'
$ws.Range('E48').Value = '  +4.27%  '

$ws.Range('D49').Value = '2.455.98'
$ws.Range('E49').Value = '  +11.42%  '

$ws.Range('E50').Value = '  +3.70%  '

$ws.Range('D51').Value = '''2.37'
$ws.Range('E51').Value = '  +15.38%  '
